$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 6685.263
$ws.Range("I28").Value = 302.93332
$ws.Range("J28").Value = 30619
$ws.Range("K28").Value = 302.93332
$ws.Range("L28").Value = 30619
$ws.Range("M28").Value = 182.06668
$ws.Range("N28").Value = -31589
$ws.Range("H69").Value = 3199.8
$ws.Range("J69").Value = 3199.8
$ws.Range("L69").Value = 9599.400000000001
$ws.Range("N69").Value = -11347.4
$ws.Range("H72").Value = 3199.8
$ws.Range("J72").Value = 3199.8
$ws.Range("L72").Value = 28798.2
$ws.Range("N72").Value = -37534.2
$ws.Range("H103").Value = 111594.336
$ws.Range("J103").Value = 560
$ws.Range("L103").Value = 1680
$ws.Range("N103").Value = -2852
$ws.Range("H129").Value = 1135.7273
$ws.Range("J129").Value = 1521.45
$ws.Range("L129").Value = 4564.35
$ws.Range("N129").Value = -14564.35
$ws.Range("H137").Value = 3528.348
$ws.Range("I137").Value = 4221.857
$ws.Range("J137").Value = 2449.5557
$ws.Range("K137").Value = 12665.571
$ws.Range("L137").Value = 7348.6671
$ws.Range("M137").Value = -10115.571
$ws.Range("N137").Value = -12448.6671
$ws.Range("H138").Value = 120683.31
$ws.Range("I138").Value = 1725.4642
$ws.Range("J138").Value = 172727.38
$ws.Range("K138").Value = 5176.392599999999
$ws.Range("L138").Value = 518182.14
$ws.Range("M138").Value = -36.39259999999922
$ws.Range("N138").Value = -528462.14

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 618506.8
$ws.Range("I32").Value = 698078.2
$ws.Range("K32").Value = 698078.2
$ws.Range("M32").Value = -697791.2
$ws.Range("H61").Value = 3300
$ws.Range("I61").Value = 2250
$ws.Range("K61").Value = 2250
$ws.Range("M61").Value = -2038
$ws.Range("H74").Value = 1329.7307
$ws.Range("I74").Value = 1009.1053
$ws.Range("J74").Value = 2200
$ws.Range("K74").Value = 1009.1053
$ws.Range("L74").Value = 2200
$ws.Range("M74").Value = -135.1053000000001
$ws.Range("N74").Value = -3948
$ws.Range("H77").Value = 1329.7307
$ws.Range("I77").Value = 1009.1053
$ws.Range("J77").Value = 2200
$ws.Range("K77").Value = 5045.5265
$ws.Range("L77").Value = 11000
$ws.Range("M77").Value = -677.5264999999999
$ws.Range("N77").Value = -19736
$ws.Range("H135").Value = 46882
$ws.Range("J135").Value = 46882
$ws.Range("L135").Value = 46882
$ws.Range("N135").Value = -57022
$ws.Range("H136").Value = 3300
$ws.Range("I136").Value = 2250
$ws.Range("K136").Value = 6750
$ws.Range("M136").Value = -4200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1288.1177
$ws.Range("I80").Value = 2760.5715
$ws.Range("J80").Value = 257.4
$ws.Range("K80").Value = 2760.5715
$ws.Range("L80").Value = 257.4
$ws.Range("M80").Value = -1762.5715
$ws.Range("N80").Value = -2253.4
$ws.Range("H83").Value = 1288.1177
$ws.Range("I83").Value = 2760.5715
$ws.Range("J83").Value = 257.4
$ws.Range("K83").Value = 13802.8575
$ws.Range("L83").Value = 1287
$ws.Range("M83").Value = -8810.8575
$ws.Range("N83").Value = -11271
$ws.Range("H107").Value = 1599.7778
$ws.Range("I107").Value = 1359.6
$ws.Range("J107").Value = 1900
$ws.Range("K107").Value = 1359.6
$ws.Range("L107").Value = 1900
$ws.Range("M107").Value = 560.4000000000001
$ws.Range("N107").Value = -5740
$ws.Range("H134").Value = 3500.2068
$ws.Range("I134").Value = 3239.3914
$ws.Range("J134").Value = 4500
$ws.Range("K134").Value = 9718.174199999999
$ws.Range("L134").Value = 13500
$ws.Range("M134").Value = -7183.174199999999
$ws.Range("N134").Value = -18570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 100000
$ws.Range("J12").Value = 100000
$ws.Range("L12").Value = 100000
$ws.Range("N12").Value = -100340
$ws.Range("H31").Value = 8943.695
$ws.Range("I31").Value = 2541.5715
$ws.Range("J31").Value = 11744.625
$ws.Range("K31").Value = 2541.5715
$ws.Range("L31").Value = 11744.625
$ws.Range("M31").Value = -2246.5715
$ws.Range("N31").Value = -12334.625
$ws.Range("H34").Value = 8943.695
$ws.Range("I34").Value = 2541.5715
$ws.Range("J34").Value = 11744.625
$ws.Range("K34").Value = 2541.5715
$ws.Range("L34").Value = 11744.625
$ws.Range("M34").Value = -2339.5715
$ws.Range("N34").Value = -12148.625
$ws.Range("H58").Value = 1188.8334
$ws.Range("I58").Value = 869.5
$ws.Range("J58").Value = 1588
$ws.Range("K58").Value = 869.5
$ws.Range("L58").Value = 1588
$ws.Range("M58").Value = -666.5
$ws.Range("N58").Value = -1994
$ws.Range("H132").Value = 11906973
$ws.Range("I132").Value = 1092.8572
$ws.Range("J132").Value = 23812852
$ws.Range("K132").Value = 3278.5716
$ws.Range("L132").Value = 71438556
$ws.Range("M132").Value = -748.5715999999998
$ws.Range("N132").Value = -71443616
$ws.Range("H134").Value = 3520.6667
$ws.Range("I134").Value = 963.5
$ws.Range("J134").Value = 8635
$ws.Range("K134").Value = 2890.5
$ws.Range("L134").Value = 25905
$ws.Range("M134").Value = -355.5
$ws.Range("N134").Value = -30975
$ws.Range("H136").Value = 1188.8334
$ws.Range("I136").Value = 869.5
$ws.Range("J136").Value = 1588
$ws.Range("K136").Value = 2608.5
$ws.Range("L136").Value = 4764
$ws.Range("M136").Value = -58.5
$ws.Range("N136").Value = -9864

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 326.72415
$ws.Range("I107").Value = 253.54546
$ws.Range("J107").Value = 556.7143
$ws.Range("K107").Value = 760.6363799999999
$ws.Range("L107").Value = 1670.1429
$ws.Range("M107").Value = 1159.36362
$ws.Range("N107").Value = -5510.1429
$ws.Range("H113").Value = 870.7544
$ws.Range("I113").Value = 669.85297
$ws.Range("J113").Value = 1167.7391
$ws.Range("K113").Value = 2009.55891
$ws.Range("L113").Value = 3503.2173
$ws.Range("M113").Value = 160.4410899999998
$ws.Range("N113").Value = -7843.2173
$ws.Range("H132").Value = 3617.3442
$ws.Range("I132").Value = 2553.2307
$ws.Range("J132").Value = 4407.8286
$ws.Range("K132").Value = 22979.0763
$ws.Range("L132").Value = 39670.4574
$ws.Range("M132").Value = -20449.0763
$ws.Range("N132").Value = -44730.4574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 4349.95
$ws.Range("J46").Value = 4349.95
$ws.Range("L46").Value = 4349.95
$ws.Range("N46").Value = -4661.95
$ws.Range("H122").Value = 4386.375
$ws.Range("I122").Value = 2835.75
$ws.Range("J122").Value = 4903.25
$ws.Range("K122").Value = 8507.25
$ws.Range("L122").Value = 14709.75
$ws.Range("M122").Value = -6057.25
$ws.Range("N122").Value = -19609.75
$ws.Range("H126").Value = 2673
$ws.Range("I126").Value = 2673
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8019
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5549
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1671.3334
$ws.Range("I16").Value = 1522.7778
$ws.Range("J16").Value = 2117
$ws.Range("K16").Value = 1522.7778
$ws.Range("L16").Value = 2117
$ws.Range("M16").Value = -1352.7778
$ws.Range("N16").Value = -2457
$ws.Range("H55").Value = 686.4828
$ws.Range("I55").Value = 255.11111
$ws.Range("J55").Value = 880.6
$ws.Range("K55").Value = 255.11111
$ws.Range("L55").Value = 880.6
$ws.Range("M55").Value = -82.11111
$ws.Range("N55").Value = -1226.6
$ws.Range("H132").Value = 4387.9565
$ws.Range("I132").Value = 3968.8
$ws.Range("K132").Value = 11906.4
$ws.Range("M132").Value = -9376.400000000001
$ws.Range("H136").Value = 10418619
$ws.Range("I136").Value = 2568
$ws.Range("J136").Value = 16668250
$ws.Range("K136").Value = 7704
$ws.Range("L136").Value = 50004750
$ws.Range("M136").Value = -5154
$ws.Range("N136").Value = -50009850

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 37055.855
$ws.Range("J123").Value = 37055.855
$ws.Range("L123").Value = 37055.855
$ws.Range("N123").Value = -46855.855
$ws.Range("H132").Value = 9807011
$ws.Range("I132").Value = 3613.875
$ws.Range("J132").Value = 18521142
$ws.Range("K132").Value = 10841.625
$ws.Range("L132").Value = 55563426
$ws.Range("M132").Value = -8311.625
$ws.Range("N132").Value = -55568486
$ws.Range("H136").Value = 3204.4211
$ws.Range("I136").Value = 2830.92
$ws.Range("K136").Value = 8492.76
$ws.Range("M136").Value = -5942.76
